# auto_UI/testcase/data/test_login.xlsx - original excel for testcase
#
# The "expected"/"actual" result columns (E:F) on the data rows had been
# filled in with a sample validation run (error text + Pass/Fail, plus a
# red/green highlight on the first data row). This restores the sheet to
# its pristine, un-run state: the E:F sample values are cleared and the
# one-off highlight on row 2 is removed so every data row shares the same
# plain formatting again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (A2:F2) was highlighted (red fill) to flag the sample run - put it
# back to no fill so it matches the rest of the table.
$ws.Range("A2:F2").Interior.Pattern = -4142

# Clear the leftover "expected"/"actual" sample output in columns E:F for
# every data row (2-5); only the static columns A:D (id/username/password/
# case) stay populated.
$ws.Range("E2:F5").ClearContents()

# Restore the plain (non-range) selection that a freshly opened sheet has.
$ws.Range("D8").Select()

Write-Host "test_login.xlsx restored to pristine testcase state"
